# Dataset Log-trasformato Scalato XGBoost scelto
# Refresh the metrics table (Random State + Training/Test RMSE/MSE/R^2 + Training Time)
# with the results produced after switching to the log-transformed, scaled dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18
$ws.Range("C2").Value = 0.09482646153591366
$ws.Range("D2").Value = 0.008992057807422113
$ws.Range("E2").Value = 0.7336261084701818
$ws.Range("F2").Value = 0.02312480534815172
$ws.Range("G2").Value = 0.0005347566223899063
$ws.Range("H2").Value = 0.9847742714255313
$ws.Range("I2").Value = 0

$ws.Range("B3").Value = 80
$ws.Range("C3").Value = 0.1195654240444575
$ws.Range("D3").Value = 0.01429589062693093
$ws.Range("E3").Value = 0.6314514225018637
$ws.Range("F3").Value = 0.03873184960659918
$ws.Range("G3").Value = 0.001500156173948217
$ws.Range("H3").Value = 0.9028935794091193
$ws.Range("I3").Value = 0.002990961074829102

$ws.Range("C4").Value = 0.18373768617052
$ws.Range("D4").Value = 0.03375953731929648
$ws.Range("F4").Value = 0.1873855451473129
$ws.Range("G4").Value = 0.03511334253015564
$ws.Range("H4").Value = [double]"-1.546049510636749e-09"
$ws.Range("I4").Value = 0.002991437911987305

$ws.Range("B5").Value = 80
$ws.Range("F5").Value = 0.006432844888492792
$ws.Range("G5").Value = [double]"4.138149335940784e-05"
$ws.Range("H5").Value = 0.9973213397587389
$ws.Range("I5").Value = 0.002347230911254883

$ws.Range("B6").Value = 83
$ws.Range("C6").Value = 0.04631584119436853
$ws.Range("D6").Value = 0.002145157145541965
$ws.Range("E6").Value = 0.9472750420931946
$ws.Range("F6").Value = 0.002916006050421347
$ws.Range("G6").Value = [double]"8.503091286093906e-06"
$ws.Range("H6").Value = 0.998845769281683
$ws.Range("I6").Value = 0.1521124839782715

$ws.Range("B7").Value = 36
$ws.Range("C7").Value = 0.1372283170240384
$ws.Range("D7").Value = 0.01883161099324999
$ws.Range("E7").Value = 0.3529241615997377
$ws.Range("F7").Value = 0.1036722936857353
$ws.Range("G7").Value = 0.01074794447806135
$ws.Range("H7").Value = 0.7619592110376481
$ws.Range("I7").Value = 0.0009980201721191406

$ws.Range("B8").Value = 61
$ws.Range("C8").Value = [double]"1.250252383024717e-07"
$ws.Range("D8").Value = [double]"1.563131021258985e-14"
$ws.Range("E8").Value = 0.9999999999995979
$ws.Range("F8").Value = 0.03445772132397083
$ws.Range("G8").Value = 0.001187334558840434
$ws.Range("H8").Value = 0.922292498210229
$ws.Range("I8").Value = 0.00197148323059082

$ws.Range("B9").Value = 72
$ws.Range("C9").Value = 0.1061764900883697
$ws.Range("D9").Value = 0.01127344704748566
$ws.Range("E9").Value = 0.6864283876845669
$ws.Range("F9").Value = 0.03742826243096131
$ws.Range("G9").Value = 0.00140087482860091
$ws.Range("H9").Value = 0.9468672020306033
$ws.Range("I9").Value = 0.01784348487854004

$ws.Range("B10").Value = 80
$ws.Range("C10").Value = 0.0007202721970399165
$ws.Range("D10").Value = [double]"5.187920378287083e-07"
$ws.Range("E10").Value = 0.9999866255224982
$ws.Range("F10").Value = 0.006526189894802719
$ws.Range("G10").Value = [double]"4.259115454302513e-05"
$ws.Range("H10").Value = 0.9972430373328257
$ws.Range("I10").Value = 0.04020166397094727
